# Generate Report for Archive
#
# This applies the localization-status report refresh:
#   - Rows for "3cc29baf-...md" and "997c5ca3-...md" swap places (row 4 <-> row 5)
#     on the Overview, zh-cn and de-de sheets.
#   - The status for "3981ea83-...md" and "997c5ca3-...md" moves from
#     "Ready for handoff" to "In Translation" on all three sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("B3").Value = "In Translation"
$ws.Range("C3").Value = "In Translation"

$ws.Range("A4").Value = "997c5ca3-bd55-40d7-8c0c-7e2238d6cc88.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"
$ws.Range("D4").Value = "2016-35-11 10:35:40"

$ws.Range("A5").Value = "3cc29baf-975a-41bd-9962-69bf3f49795b.md"
$ws.Range("D5").Value = "2016-34-11 10:34:51"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("C3").Value = "In Translation"

$ws.Range("A4").Value = "997c5ca3-bd55-40d7-8c0c-7e2238d6cc88.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("D4").Value = "997c5ca3-bd55-40d7-8c0c-7e2238d6cc88.ae441ff79468e811d444e1fe999aa105c4187e40.zh-cn.xlf"
$ws.Range("E4").Value = "2016-03-11 10:35:37"

$ws.Range("A5").Value = "3cc29baf-975a-41bd-9962-69bf3f49795b.md"
$ws.Range("D5").Value = "3cc29baf-975a-41bd-9962-69bf3f49795b.57055bd7c49ef0346b01e57de0bf2059c29d2b85.zh-cn.xlf"
$ws.Range("E5").Value = "2016-03-11 10:34:46"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("C3").Value = "In Translation"

$ws.Range("A4").Value = "997c5ca3-bd55-40d7-8c0c-7e2238d6cc88.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("D4").Value = "997c5ca3-bd55-40d7-8c0c-7e2238d6cc88.ae441ff79468e811d444e1fe999aa105c4187e40.de-de.xlf"
$ws.Range("E4").Value = "2016-03-11 10:35:40"

$ws.Range("A5").Value = "3cc29baf-975a-41bd-9962-69bf3f49795b.md"
$ws.Range("D5").Value = "3cc29baf-975a-41bd-9962-69bf3f49795b.57055bd7c49ef0346b01e57de0bf2059c29d2b85.de-de.xlf"
$ws.Range("E5").Value = "2016-03-11 10:34:51"
